$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at the top of the Cebolla data block (row 296),
# pushing the existing rows (old 296..351) down to 298..353.
$ws.Rows.Item(296).Resize(2).Insert()

# Row 296: new "1a nueva(o)" entry for the latest reporting date.
$ws.Range("A296").Value = 7
$ws.Range("B296").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C296").Value = "Ñuble"
$ws.Range("D296").Value = 44476
$ws.Range("E296").Value = 16
$ws.Range("F296").Value = 100112004
$ws.Range("G296").Value = "Cebolla"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "1a nueva(o)"
$ws.Range("J296").Value = 18000
$ws.Range("K296").Value = 1300
$ws.Range("L296").Value = 1400
$ws.Range("M296").Value = 1350
$ws.Range("N296").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O296").Value = "Región de O'Higgins"
$ws.Range("P296").Value = 135
$ws.Range("Q296").Value = 10
$ws.Range("R296").Value = "Hortaliza"

# Row 297: new "2a nueva(o)" entry for the latest reporting date.
$ws.Range("A297").Value = 7
$ws.Range("B297").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C297").Value = "Ñuble"
$ws.Range("D297").Value = 44476
$ws.Range("E297").Value = 16
$ws.Range("F297").Value = 100112004
$ws.Range("G297").Value = "Cebolla"
$ws.Range("H297").Value = "Sin especificar"
$ws.Range("I297").Value = "2a nueva(o)"
$ws.Range("J297").Value = 12000
$ws.Range("K297").Value = 1100
$ws.Range("L297").Value = 1200
$ws.Range("M297").Value = 1150
$ws.Range("N297").Value = "$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O297").Value = "Región de O'Higgins"
$ws.Range("P297").Value = 115
$ws.Range("Q297").Value = 10
$ws.Range("R297").Value = "Hortaliza"
